# "wrapping up test file audit"
#
# The "optimization_parameters" sheet had a stray leftover row (row 16:
# label "Sheet" with values 3 / 4) left over from an earlier edit that was
# never cleaned up. Delete that whole row -- the rows below it (most
# notably the "simulation_timepoints" row) shift up to fill the gap.
#
# Deleting that row also removes the only cell that used the "Sheet"
# shared string and the only cell that used the integer (numFmtId=1)
# number format, so both get dropped from the workbook's shared tables
# when it is saved -- that is an automatic side effect of this edit, not
# a separate step.

$wb = $excel.ActiveWorkbook

$paramSheet = $wb.Worksheets.Item("optimization_parameters")
$paramSheet.Rows.Item(16).Delete()

# Leave the selection sitting on the row that just shifted up into row 16
# (now "simulation_timepoints"), matching what Excel does after an
# Edit > Delete > Entire Row operation.
$paramSheet.Range("A16:XFD16").Select()

# The workbook was left with "degradation_rates" as the active/selected
# sheet rather than "optimization_parameters".
$degSheet = $wb.Worksheets.Item("degradation_rates")
$degSheet.Activate()
$degSheet.Range("B1").Select()
